# Adds a new client row ("FRANK FERRETERIA FRANKFERRE CIA.") for
# ILLER LOPEZ ROBERTO FERNANDO on both worksheets, right before the
# existing "VIEJO RIVAS MAYRA ANABELLE" row, renames the prior B7
# client to "CORPORACION AREVALO-YUMBLA E HIJOS", widens column B by
# one unit, and bumps the trailing totals row's "0 de N" labels
# (where present) to reflect the new row count.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count
    $lastCol = $usedRange.Columns.Count

    # Widen column B from 35 to 36. The COM ColumnWidth property adds a
    # constant ~0.8333 padding offset relative to the raw OOXML <col
    # width> units, so subtract it to land exactly on 36 in the saved
    # file (matches how the original 35 round-trips to a clean integer).
    $ws.Columns.Item(2).ColumnWidth = 35.166666666666664

    # B7 used to hold "FRANK FERRETERIA FRANKFERRE CIA."; it now holds a
    # different client, and that ferretería entry becomes its own new row.
    $ws.Cells.Item(7, 2).Value = "CORPORACION AREVALO-YUMBLA E HIJOS"

    # Insert a new row above the old row 8 ("VIEJO RIVAS ..."), which
    # pushes it (and the totals row after it) one row down.
    $ws.Rows.Item(8).Insert()

    $ws.Cells.Item(8, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
    $ws.Cells.Item(8, 2).Value = "FRANK FERRETERIA FRANKFERRE CIA."
    for ($col = 3; $col -le $lastCol; $col++) {
        $ws.Cells.Item(8, $col).Value = 0
    }

    # The totals row, now shifted down by one, needs its "0 de 7" style
    # labels (if any) updated to "0 de 8" to reflect the extra data row.
    $newTotalsRow = $lastRow + 1
    for ($col = 3; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($newTotalsRow, $col)
        if ($cell.Value2 -eq "0 de 7") {
            $cell.Value = "0 de 8"
        }
    }
}
